$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row 1 / A1 "ElecTweets" already present from before)
$ws.Range("A2").Value = "Elections letweetlatives "
$ws.Range("B2").Value = ":p"
$ws.Range("A3").Value = "Future is tweet"
$ws.Range("B3").Value = "Jeu de mot avec `"futuristic`"… mais pas évident !"

# Column widths (closest representable values to the target 23.44140625 / 40.21875
# given this engine's internal character-width quantization)
$ws.Columns.Item(1).ColumnWidth = 22.666666666666668
$ws.Columns.Item(2).ColumnWidth = 39.333333333333336

# Update the active selection to match the saved view state
$ws.Range("B8").Select() | Out-Null
